# C4-PowerPoint.pptx theme update
#
# The authored change swaps the presentation's colour theme from the
# custom "Integral" palette over to the built-in "Office Theme" palette
# (the two theme parts already shared identical font/format schemes, so
# only the colour scheme actually changes visually). We reproduce that
# by rewriting the twelve theme colour slots on the active theme through
# the documented COM surface:
#   ActivePresentation.SlideMaster.Theme.ThemeColorScheme.Colors(i).RGB
#
# RGB() packs as 0x00BBGGRR (classic VBA colour order), so convert each
# target "RRGGBB" hex swatch accordingly before assigning it.

function Get-RgbFromHex($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# clrScheme slot order exposed via Colors(1..12):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeTheme = @{
    1  = "000000"  # dk1      (unchanged)
    2  = "FFFFFF"  # lt1      (unchanged)
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

for ($i = 1; $i -le 12; $i++) {
    $colors.Colors($i).RGB = Get-RgbFromHex $officeTheme[$i]
}
